$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 257, pushing the existing 257.. rows down to 259..
$ws.Rows.Item(257).Insert()
$ws.Rows.Item(257).Insert()

# New row 257 - updated "Patagonia / 1a (guarda)" record (date + price changes)
$ws.Cells.Item(257, 1).Value = 4
$ws.Cells.Item(257, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(257, 3).Value = "Los Lagos"
$ws.Cells.Item(257, 4).Value = 44509
$ws.Cells.Item(257, 5).Value = 10
$ws.Cells.Item(257, 6).Value = 100114001
$ws.Cells.Item(257, 7).Value = "Papa"
$ws.Cells.Item(257, 8).Value = "Patagonia"
$ws.Cells.Item(257, 9).Value = "1a (guarda)"
$ws.Cells.Item(257, 10).Value = 300
$ws.Cells.Item(257, 11).Value = 8000
$ws.Cells.Item(257, 12).Value = 9000
$ws.Cells.Item(257, 13).Value = 8500
$ws.Cells.Item(257, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(257, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(257, 16).Value = 340
$ws.Cells.Item(257, 17).Value = 25
$ws.Cells.Item(257, 18).Value = "Hortaliza"

# New row 258 - brand-new "Pehuenche / 1a nueva(o)" record
$ws.Cells.Item(258, 1).Value = 4
$ws.Cells.Item(258, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(258, 3).Value = "Los Lagos"
$ws.Cells.Item(258, 4).Value = 44509
$ws.Cells.Item(258, 5).Value = 10
$ws.Cells.Item(258, 6).Value = 100114001
$ws.Cells.Item(258, 7).Value = "Papa"
$ws.Cells.Item(258, 8).Value = "Pehuenche"
$ws.Cells.Item(258, 9).Value = "1a nueva(o)"
$ws.Cells.Item(258, 10).Value = 300
$ws.Cells.Item(258, 11).Value = 16000
$ws.Cells.Item(258, 12).Value = 16000
$ws.Cells.Item(258, 13).Value = 16000
$ws.Cells.Item(258, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(258, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(258, 16).Value = 640
$ws.Cells.Item(258, 17).Value = 25
$ws.Cells.Item(258, 18).Value = "Hortaliza"
